# Delete the rows for "Resolving-Mac" sending cluster (former rows 11-13)
$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Rows("11:13").Delete()

# Update remaining data rows (2-10) with refreshed TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.835339
$ws.Range("H2").Value = 11.506017
$ws.Range("I2").Value = 0.4358468215464834
$ws.Range("J2").Value = 0.4358468215464834
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07259900000000001
$ws.Range("N2").Value = 0.217797
$ws.Range("O2").Value = 0.0162094769588191
$ws.Range("P2").Value = 0.0162094769588191
$ws.Range("Q2").Value = 0.278441776061
$ws.Range("R2").Value = 2.505975984549
$ws.Range("S2").Value = 0.007064849011432263
$ws.Range("T2").Value = 0.007064849011432263

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.835339
$ws.Range("H3").Value = 11.506017
$ws.Range("I3").Value = 0.4358468215464834
$ws.Range("J3").Value = 0.4358468215464834
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.187332
$ws.Range("N3").Value = 12.561996
$ws.Range("O3").Value = 0.9349228167457665
$ws.Range("P3").Value = 0.9349228167457664
$ws.Range("Q3").Value = 16.059837725548
$ws.Range("R3").Value = 144.538539529932
$ws.Range("S3").Value = 0.4074831380699277
$ws.Range("T3").Value = 0.4074831380699276

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.835339
$ws.Range("H4").Value = 11.506017
$ws.Range("I4").Value = 0.4358468215464834
$ws.Range("J4").Value = 0.4358468215464834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2188686666666667
$ws.Range("N4").Value = 0.656606
$ws.Range("O4").Value = 0.04886770629541442
$ws.Range("P4").Value = 0.04886770629541441
$ws.Range("Q4").Value = 0.8394355331446667
$ws.Range("R4").Value = 7.554919798302
$ws.Range("S4").Value = 0.02129883446512345
$ws.Range("T4").Value = 0.02129883446512345

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.744136
$ws.Range("H5").Value = 5.232408
$ws.Range("I5").Value = 0.1982031137129723
$ws.Range("J5").Value = 0.1982031137129723
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07259900000000001
$ws.Range("N5").Value = 0.217797
$ws.Range("O5").Value = 0.0162094769588191
$ws.Range("P5").Value = 0.0162094769588191
$ws.Range("Q5").Value = 0.126622529464
$ws.Range("R5").Value = 1.139602765176
$ws.Range("S5").Value = 0.003212768804896627
$ws.Range("T5").Value = 0.003212768804896627

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.744136
$ws.Range("H6").Value = 5.232408
$ws.Range("I6").Value = 0.1982031137129723
$ws.Range("J6").Value = 0.1982031137129723
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.187332
$ws.Range("N6").Value = 12.561996
$ws.Range("O6").Value = 0.9349228167457665
$ws.Range("P6").Value = 0.9349228167457664
$ws.Range("Q6").Value = 7.303276485152002
$ws.Range("R6").Value = 65.729488366368
$ws.Range("S6").Value = 0.1853046133603135
$ws.Range("T6").Value = 0.1853046133603135

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.744136
$ws.Range("H7").Value = 5.232408
$ws.Range("I7").Value = 0.1982031137129723
$ws.Range("J7").Value = 0.1982031137129723
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2188686666666667
$ws.Range("N7").Value = 0.656606
$ws.Range("O7").Value = 0.04886770629541442
$ws.Range("P7").Value = 0.04886770629541441
$ws.Range("Q7").Value = 0.3817367208053334
$ws.Range("R7").Value = 3.435630487248
$ws.Range("S7").Value = 0.009685731547762156
$ws.Range("T7").Value = 0.009685731547762155

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Bmpr1b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.220265666666667
$ws.Range("H8").Value = 9.660797
$ws.Range("I8").Value = 0.3659500647405443
$ws.Range("J8").Value = 0.3659500647405443
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07259900000000001
$ws.Range("N8").Value = 0.217797
$ws.Range("O8").Value = 0.0162094769588191
$ws.Range("P8").Value = 0.0162094769588191
$ws.Range("Q8").Value = 0.2337880671343334
$ws.Range("R8").Value = 2.104092604209
$ws.Range("S8").Value = 0.00593185914249021
$ws.Range("T8").Value = 0.00593185914249021

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Bmpr1b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.220265666666667
$ws.Range("H9").Value = 9.660797
$ws.Range("I9").Value = 0.3659500647405443
$ws.Range("J9").Value = 0.3659500647405443
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.187332
$ws.Range("N9").Value = 12.561996
$ws.Range("O9").Value = 0.9349228167457665
$ws.Range("P9").Value = 0.9349228167457664
$ws.Range("Q9").Value = 13.48432147453467
$ws.Range("R9").Value = 121.358893270812
$ws.Range("S9").Value = 0.3421350653155252
$ws.Range("T9").Value = 0.3421350653155252

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Bmpr1b"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.220265666666667
$ws.Range("H10").Value = 9.660797
$ws.Range("I10").Value = 0.3659500647405443
$ws.Range("J10").Value = 0.3659500647405443
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2188686666666667
$ws.Range("N10").Value = 0.656606
$ws.Range("O10").Value = 0.04886770629541442
$ws.Range("P10").Value = 0.04886770629541441
$ws.Range("Q10").Value = 0.7048152527757778
$ws.Range("R10").Value = 6.343337274982001
$ws.Range("S10").Value = 0.01788314028252881
$ws.Range("T10").Value = 0.0178831402825288
